# Apply the "output generated at 456a3b4" refresh to 北京-漫展信息.xlsx
# Sheets: 1=展览 (Exhibition), 2=演出 (Performance), 3=本地生活 (Local life), 4=全部类型 (All types)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet 1 (展览): 想去人数 (F column) refreshed counts ---
$ws1.Range("F6").Value  = 224
$ws1.Range("F7").Value  = 205
$ws1.Range("F8").Value  = 234
$ws1.Range("F9").Value  = 2821
$ws1.Range("F12").Value = 2162
$ws1.Range("F13").Value = 245
$ws1.Range("F17").Value = 2509
$ws1.Range("F19").Value = 1246
$ws1.Range("F20").Value = 4518
$ws1.Range("F22").Value = 4744
$ws1.Range("F23").Value = 1283
$ws1.Range("F24").Value = 2766
$ws1.Range("F25").Value = 3171
$ws1.Range("F26").Value = 136
$ws1.Range("F27").Value = 1478
$ws1.Range("F29").Value = 816
$ws1.Range("F30").Value = 84
$ws1.Range("F31").Value = 236
$ws1.Range("F32").Value = 835
$ws1.Range("F33").Value = 1441
$ws1.Range("F34").Value = 108
$ws1.Range("F35").Value = 225
$ws1.Range("F36").Value = 572
$ws1.Range("F37").Value = 149
$ws1.Range("F38").Value = 276
$ws1.Range("F39").Value = 348

# --- Sheet 2 (演出): 想去人数 refreshed count ---
$ws2.Range("F7").Value = 46

# --- Sheet 4 (全部类型): row 6 replaced with a different event entirely ---
$ws4.Range("C6").Value = "北京·【母亲节 5折 特惠】“她”的协奏曲——《致爱丽丝》《欢乐颂》《假如爱有天意》母亲节专场烛光音乐会"
$ws4.Range("D6").Value = "朝阳北路常营陆港城20号院1号楼 常营·爱乐汇艺术空间(长楹天街店)"
$ws4.Range("E6").Value = "2024.05.12 15:30-05.12 17:00"
$ws4.Range("F6").Value = 1
$ws4.Range("G6").Value = 80
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=84818"
$ws4.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202404/jnV1IFn41714013668157.jpeg"

# --- Sheet 4 (全部类型): 想去人数 refreshed counts ---
$ws4.Range("F8").Value  = 224
$ws4.Range("F9").Value  = 205
$ws4.Range("F10").Value = 46
$ws4.Range("F11").Value = 234
$ws4.Range("F12").Value = 2821
$ws4.Range("F14").Value = 117
$ws4.Range("F15").Value = 2162
$ws4.Range("F16").Value = 245
$ws4.Range("F21").Value = 2509
$ws4.Range("F22").Value = 1246
$ws4.Range("F26").Value = 4518
$ws4.Range("F28").Value = 4744
$ws4.Range("F29").Value = 1283
$ws4.Range("F30").Value = 2766
$ws4.Range("F31").Value = 3171
$ws4.Range("F32").Value = 136
$ws4.Range("F35").Value = 1478
$ws4.Range("F38").Value = 816
$ws4.Range("F39").Value = 84
$ws4.Range("F40").Value = 236
$ws4.Range("F41").Value = 835
$ws4.Range("F43").Value = 1441
$ws4.Range("F44").Value = 108
$ws4.Range("F45").Value = 225
$ws4.Range("F46").Value = 572
$ws4.Range("F47").Value = 149
$ws4.Range("F48").Value = 276
$ws4.Range("F49").Value = 348

Write-Output "Applied all updates"
